$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.809.37"
$ws.Range("E2").Value = "  +1.72%  "
$ws.Range("D3").Value = "3.736.16"
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.15"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.45"
$ws.Range("E6").Value = "  -5.46%  "
$ws.Range("D7").Value = "3.734.76"
$ws.Range("E7").Value = "  -1.96%  "
$ws.Range("E9").Value = "  +1.95%  "
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.36"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.99"
$ws.Range("E13").Value = "  -0.89%  "
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "4.363.49"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "3.735.03"
$ws.Range("E16").Value = "  -1.93%  "
$ws.Range("D17").Value = "68.787.95"
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.24"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.32"
$ws.Range("E20").Value = "  +4.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.52"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("E22").Value = "  +10.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.89"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000141"
$ws.Range("E26").Value = "  -7.57%  "
$ws.Range("E27").Value = "  +0.63%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.44"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.90"
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.63"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").Value = "3.882.65"
$ws.Range("E34").Value = "  -1.72%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "3.669.95"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("E38").Value = "  +0.39%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  -2.08%  "
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "433.32"
$ws.Range("E42").Value = "  -4.01%  "
$ws.Range("E44").Value = "  -1.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.86"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.59"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.00"
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").Value = "2.743.87"
$ws.Range("E51").Value = "  -3.56%  "
